$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D4").Value = 1.0407973260155801
$ws.Range("D5").Value = 1.0476333876758299
$ws.Range("D6").Value = 1.04643379828476
$ws.Range("D7").Value = 1.03057960659685
$ws.Range("D12").Value = 1.0407973260155801
$ws.Range("D13").Value = 1.0096340779125801
$ws.Range("D14").Value = 0.97271823593722395
$ws.Range("D15").Value = 0.86955621894563395
$ws.Range("D20").Value = 1.09147654328549
$ws.Range("D21").Value = 1.0607662171155801
$ws.Range("D22").Value = 1.0211354838968201
$ws.Range("D23").Value = 0.91204204096250496
$ws.Range("D28").Value = 1.0407907558492699
$ws.Range("D29").Value = 1.0076701024091199
$ws.Range("D30").Value = 0.96931667960613599
$ws.Range("D31").Value = 0.86230827052423598
$ws.Range("D36").Value = 1.0535659467793399
$ws.Range("D37").Value = 1.0212884963159301
$ws.Range("D38").Value = 0.98221107696883803
$ws.Range("D39").Value = 0.87381315481220201
$ws.Range("D44").Value = 1.05345859447816
$ws.Range("D45").Value = 1.0211204549644699
$ws.Range("D46").Value = 0.98208786760977596
$ws.Range("D47").Value = 0.87365500225397696
$ws.Range("D52").Value = 1.03910025720013
$ws.Range("D53").Value = 1.00597506091671
$ws.Range("D54").Value = 0.96762461749252704
$ws.Range("D55").Value = 0.86082265813353398
$ws.Range("D60").Value = 1.0562373267640499
$ws.Range("D61").Value = 1.0559532854573701
$ws.Range("D62").Value = 1.054412563718
$ws.Range("D63").Value = 1.0382703098046899
$ws.Range("D68").Value = 1.0562373267640499
$ws.Range("D69").Value = 1.02259103051849
$ws.Range("D70").Value = 0.98218344736982199
$ws.Range("D71").Value = 0.87735819198776099
$ws.Range("D76").Value = 1.0952905556632699
$ws.Range("D77").Value = 1.0621260212758299
$ws.Range("D78").Value = 1.0201085030949999
$ws.Range("D79").Value = 0.90911373637535098
$ws.Range("D84").Value = 1.0573160138432001
$ws.Range("D85").Value = 1.02120128529842
$ws.Range("D86").Value = 0.97951079936762697
$ws.Range("D87").Value = 0.87128586905044003
$ws.Range("D92").Value = 1.0704249800211301
$ws.Range("D93").Value = 1.03482545825172
$ws.Range("D94").Value = 0.99299438925022898
$ws.Range("D95").Value = 0.88243913296169096
$ws.Range("D100").Value = 1.06964750903626
$ws.Range("D101").Value = 1.0339785418535701
$ws.Range("D102").Value = 0.99208708815006497
$ws.Range("D103").Value = 0.88191103833325202
$ws.Range("D108").Value = 1.05504445321702
$ws.Range("D109").Value = 1.0185519007421699
$ws.Range("D110").Value = 0.97705268484111596
$ws.Range("D111").Value = 0.868683125940095
